$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - row 2 and row 5 "想去人数" (want-to-go count) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 294
$ws1.Range("F5").Value = 269

# Sheet "全部类型" (All types) mirrors the same data - apply identical updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 294
$ws4.Range("F5").Value = 269
